$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.526521720180534
$ws.Range("B2").Value = -1.356031934776401

$ws.Range("A3").Value = -0.5259137153876586
$ws.Range("B3").Value = -0.6159626380292385

$ws.Range("A4").Value = -0.9828431296275644
$ws.Range("B4").Value = -0.7840647253875993

$ws.Range("A5").Value = -0.7149994499608116
$ws.Range("B5").Value = -0.6303839189168619

$ws.Range("A6").Value = 0.807323655534443
$ws.Range("B6").Value = 0.5862204236525466

$ws.Range("A7").Value = -0.1198885221354363
$ws.Range("B7").Value = -0.00001969752515742964
